# Sample Project / Main.xlsx — admin SAVE
# The "Rules" sheet cell B11 (row 11, the R40 rule-id cell) is changed
# from the text "R40" to the text "1".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces Excel to store the value as literal text
# (not the number 1), matching the shared-string cell the workbook
# diff shows ("R40" -> "1").
$ws.Range("B11").Value = "'1"
